$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 764
$ws1.Range("F4").Value = 1527
$ws1.Range("F6").Value = 99
$ws1.Range("F8").Value = 6295
$ws1.Range("F12").Value = 5306
$ws1.Range("G12").Value = 39.9
$ws1.Range("F13").Value = 31
$ws1.Range("F18").Value = 366
$ws1.Range("F19").Value = 73
$ws1.Range("F23").Value = 3789

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 89

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 89
$ws4.Range("F4").Value = 764
$ws4.Range("F5").Value = 1527
$ws4.Range("F7").Value = 99
$ws4.Range("F9").Value = 6295
$ws4.Range("F13").Value = 5306
$ws4.Range("G13").Value = 39.9
$ws4.Range("F14").Value = 31
$ws4.Range("F19").Value = 366
$ws4.Range("F20").Value = 73
$ws4.Range("F24").Value = 3789
